# Auto-generated update of cryptos list values (Sat Jan  6 23:50:42 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.317.22"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.240.77"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'307.28"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "'93.70"
$ws.Range("E6").Value = "  -6.38%  "
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").Value = "'34.50"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("D11").Value = "'0.0807"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "'7.15"
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "2.339.43"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "'0.829"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "'13.44"
$ws.Range("E16").Value = "  -4.42%  "
$ws.Range("D17").Value = "44.054.53"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "'12.02"
$ws.Range("E20").Value = "  -9.08%  "
$ws.Range("D21").Value = "'65.63"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'237.14"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "'2.93"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "'39.80"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'5.90"
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.03"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "'151.79"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").Value = "'0.0795"
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("E34").Value = "  -12.86%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("E37").Value = "  -8.97%  "
$ws.Range("D38").Value = "'3.48"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").Value = "'3.79"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").Value = "'14.27"
$ws.Range("E40").Value = "  -8.06%  "
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "1.700.37"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'82.76"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'99.38"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'4.92"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "'54.67"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "'67.29"
$ws.Range("E51").Value = "  -6.16%  "
